# Edit script: update Patient_Data sheet values (columns A & B) for rows 1-34,
# removing the previous trailing blank placeholder rows (29-31) and replacing
# them with real data, then growing the used range to A1:B34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 1..34, columns A (index 1) and B (index 2)
$data = @(
    @(0.731451117032177, 211.0819327731),
    @(2.80462184873949, 211.0819327731),
    @(4.75584136093461, 211.0819327731),
    @(5.97535355605656, 211.0819327731),
    @(6.95096331215413, 211.0819327731),
    @(7.76798524287763, 30252.7983193277),
    @(8.75512400081983, 20799.1113445378),
    @(9.81989137118261, 47689.5987394958),
    @(10.7890961262553, 52941.6470588235),
    @(11.7554826808772, 60504.5966386554),
    @(12.8456138552982, 66596.9726890756),
    @(13.8099508095921, 75840.5777310924),
    @(14.5324349251895, 83403.5273109243),
    @(15.993031358885, 85714.4285714285),
    @(16.8472022955523, 85294.2647058823),
    @(18.0695326911252, 82983.3634453781),
    @(19.6584853453576, 80042.2163865546),
    @(21.0058413609346, 75210.3319327731),
    @(21.8633428981348, 72059.1029411764),
    @(23.0820864931338, 72689.3487394958),
    @(23.941637630662, 67857.4642857142),
    @(24.9218589874974, 64075.9894957983),
    @(27.2458495593359, 58403.7773109243),
    @(28.1015576962492, 56723.1218487395),
    @(28.8376204140192, 53151.7289915966),
    @(30.0591822094691, 51471.0735294117),
    @(31.1585365853658, 50000.5),
    @(32.2558413609346, 50210.5819327731),
    @(34.4599303135888, 42857.7142857142),
    @(35.0673806107809, 44748.4516806722),
    @(36.2943226070916, 38656.0756302521),
    @(37.1426009428161, 43067.7962184874),
    @(37.8773826603812, 40546.81302521),
    @(39.0974072555851, 40126.6491596638)
)

# Rows 29-31 previously existed only as empty placeholder rows with a custom
# (non-default) row height of 13.8. Reset them to the default row height so
# that, once populated, they look identical to the other data rows (ht="15",
# customHeight="false") instead of carrying over the old custom height.
$ws.Rows.Item(29).AutoFit()
$ws.Rows.Item(30).AutoFit()
$ws.Rows.Item(31).AutoFit()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Reset the selection to match the new used range A1:B34, active cell A1
$ws.Range("A1:B34").Select()
